function Set-TextValue($cell, [string]$text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") '76.107.50'
$ws.Range("E2").Value = '  +1.68%  '

Set-TextValue $ws.Range("D3") '2.943.39'
$ws.Range("E3").Value = '  +4.65%  '

$ws.Range("E4").Value = '  +0.02%  '

Set-TextValue $ws.Range("D5") '203.94'
$ws.Range("E5").Value = '  +8.64%  '

Set-TextValue $ws.Range("D6") '600.60'
$ws.Range("E6").Value = '  +1.52%  '

$ws.Range("E7").Value = '  -0.04%  '

Set-TextValue $ws.Range("D8") '0.555'
$ws.Range("E8").Value = '  +1.78%  '

Set-TextValue $ws.Range("D9") '0.198'
$ws.Range("E9").Value = '  +5.21%  '

Set-TextValue $ws.Range("D10") '2.938.85'
$ws.Range("E10").Value = '  +4.52%  '

Set-TextValue $ws.Range("D11") '0.447'
$ws.Range("E11").Value = '  +19.31%  '

Set-TextValue $ws.Range("D12") '0.161'
$ws.Range("E12").Value = '  +0.91%  '

$ws.Range("E13").Value = '  +2.17%  '

Set-TextValue $ws.Range("D14") '3.482.53'
$ws.Range("E14").Value = '  +4.70%  '

Set-TextValue $ws.Range("D15") '28.40'
$ws.Range("E15").Value = '  +5.87%  '

Set-TextValue $ws.Range("D16") '76.104.34'
$ws.Range("E16").Value = '  +1.84%  '

$ws.Range("E17").Value = '  +2.64%  '

Set-TextValue $ws.Range("D18") '2.943.09'
$ws.Range("E18").Value = '  +4.70%  '

Set-TextValue $ws.Range("D19") '12.98'
$ws.Range("E19").Value = '  +5.95%  '

Set-TextValue $ws.Range("D20") '8.84'
$ws.Range("E20").Value = '  -0.26%  '

Set-TextValue $ws.Range("D21") '374.79'
$ws.Range("E21").Value = '  -0.39%  '

Set-TextValue $ws.Range("D22") '2.34'
$ws.Range("E22").Value = '  +3.49%  '

Set-TextValue $ws.Range("D23") '4.36'
$ws.Range("E23").Value = '  +6.62%  '

Set-TextValue $ws.Range("D24") '71.85'
$ws.Range("E24").Value = '  +1.75%  '

$ws.Range("E25").Value = '  -0.06%  '

Set-TextValue $ws.Range("D26") '4.36'
$ws.Range("E26").Value = '  +5.39%  '

$ws.Range("E27").Value = '  +4.38%  '

Set-TextValue $ws.Range("D28") '9.73'
$ws.Range("E28").Value = '  +0.72%  '

$ws.Range("E29").Value = '  +7.26%  '

$ws.Range("E30").Value = '  +0.10%  '

$ws.Range("E31").Value = '  +1.13%  '

Set-TextValue $ws.Range("D32") '7.94'
$ws.Range("E32").Value = '  +4.11%  '

Set-TextValue $ws.Range("D33") '504.57'
$ws.Range("E33").Value = '  -1.02%  '

Set-TextValue $ws.Range("D34") '1.86'
$ws.Range("E34").Value = '  +4.24%  '

$ws.Range("E35").Value = '  -0.06%  '

Set-TextValue $ws.Range("D36") '20.37'
$ws.Range("E36").Value = '  +2.77%  '

Set-TextValue $ws.Range("D37") '163.89'
$ws.Range("E37").Value = '  -0.42%  '

Set-TextValue $ws.Range("D38") '0.109'
$ws.Range("E38").Value = '  +25.98%  '

Set-TextValue $ws.Range("D39") '19.65'
$ws.Range("E39").Value = '  +1.51%  '

$ws.Range("E40").Value = '  +9.79%  '

Set-TextValue $ws.Range("D41") '0.113'
$ws.Range("E41").Value = '  -3.38%  '

Set-TextValue $ws.Range("D44") '5.03'
$ws.Range("E44").Value = '  +1.12%  '

Set-TextValue $ws.Range("D45") '1.68'
$ws.Range("E45").Value = '  +0.98%  '

Set-TextValue $ws.Range("D46") '40.23'
$ws.Range("E46").Value = '  +0.65%  '

$ws.Range("E47").Value = '  +0.18%  '

Set-TextValue $ws.Range("D48") '2.36'
$ws.Range("E48").Value = '  +2.70%  '

$ws.Range("E49").Value = '  +1.85%  '

Set-TextValue $ws.Range("D50") '3.84'
$ws.Range("E50").Value = '  +3.64%  '

# Row 42: Aave -> USDe
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D42") '1.00'
$ws.Range("E42").Value = '  +0.07%  '

# Row 43: USDe -> Aave
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D43") '181.36'
$ws.Range("E43").Value = '  -0.37%  '

# Row 51: Mantle -> InjectiveProtocol
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D51") '22.77'
$ws.Range("E51").Value = '  +9.32%  '
